$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("保險")

# --- Extend formatting for the new columns (E:K) by copying from column D ---
# Header row (bold + border style, same as B1:D1)
$ws.Range("D1").Copy()
$ws.Range("E1:K1").PasteSpecial(-4122)

# Data rows (plain style, same as B2:D4)
$ws.Range("D2").Copy()
$ws.Range("E2:K2").PasteSpecial(-4122)

$ws.Range("D3").Copy()
$ws.Range("E3:K3").PasteSpecial(-4122)

$ws.Range("D4").Copy()
$ws.Range("E4:K4").PasteSpecial(-4122)

# The "date" column (G) holds the date as literal text (e.g. "2012-04-30"),
# not an actual Excel date serial, so force a Text format before writing it.
$ws.Range("G2:G4").NumberFormat = "@"

# --- Header row: real column labels instead of stray data values ---
$ws.Range("B1").Value = "company"
$ws.Range("C1").Value = "name"
$ws.Range("D1").Value = "owner"
$ws.Range("E1").Value = "property_category"
$ws.Range("F1").Value = "category"
$ws.Range("G1").Value = "date"
$ws.Range("H1").Value = "legislator_name"
$ws.Range("I1").Value = "legislator_id"
$ws.Range("J1").Value = "source_file"
$ws.Range("K1").Value = "index"

# --- Row 2 (index 103): 國華人壽 / 保障型醫療險 ---
$ws.Range("B2").Value = "國華人壽"
$ws.Range("C2").Value = "保障型醫療險"
$ws.Range("D2").Value = "謝逸萍"
$ws.Range("E2").Value = "insurance"
$ws.Range("F2").Value = "normal"
$ws.Range("G2").Value = "2012-04-30"
$ws.Range("H2").Value = "王廷升"
$ws.Range("I2").Value = 1727
$ws.Range("J2").Value = "tmpc32d1"
$ws.Range("K2").Value = 103

# --- Row 3 (index 104): 富邦人壽 / 新終身壽險 ---
$ws.Range("B3").Value = "富邦人壽"
$ws.Range("C3").Value = "新終身壽險"
$ws.Range("D3").Value = "王廷升"
$ws.Range("E3").Value = "insurance"
$ws.Range("F3").Value = "normal"
$ws.Range("G3").Value = "2012-04-30"
$ws.Range("H3").Value = "王廷升"
$ws.Range("I3").Value = 1727
$ws.Range("J3").Value = "tmpc32d1"
$ws.Range("K3").Value = 104

# --- Row 4 (index 105): 富邦人壽 / 新終身壽險 ---
$ws.Range("B4").Value = "富邦人壽"
$ws.Range("C4").Value = "新終身壽險"
$ws.Range("D4").Value = "謝逸萍"
$ws.Range("E4").Value = "insurance"
$ws.Range("F4").Value = "normal"
$ws.Range("G4").Value = "2012-04-30"
$ws.Range("H4").Value = "王廷升"
$ws.Range("I4").Value = 1727
$ws.Range("J4").Value = "tmpc32d1"
$ws.Range("K4").Value = 105
